$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 15 de Septiembre de 2020 a las 10:21"

# Peru (row 7)
$ws.Range("B7").Value = 1073849
$ws.Range("C7").Value = 5529
$ws.Range("D7").Value = 884305
$ws.Range("E7").Value = 170759
$ws.Range("G7").Value = 150
$ws.Range("H7").Value = 18785

# Singapur (row 55)
$ws.Range("B55").Value = 57488
$ws.Range("C55").Value = 34
$ws.Range("E55").Value = 659

# Rows 89/90: swap Zambia/Croacia labels and give Croacia (row 89) fresh data,
# while Zambia (row 90) takes over the previous Croacia row's old Zambia figures.
$ws.Range("A89").Value = "Croacia"
$ws.Range("B89").Value = 13749
$ws.Range("C89").Value = 151
$ws.Range("D89").Value = 11412
$ws.Range("E89").Value = 2107
$ws.Range("G89").Value = 3
$ws.Range("H89").Value = 230

$ws.Range("A90").Value = "Zambia"
$ws.Range("B90").Value = 13720
$ws.Range("C90").Value = 0
$ws.Range("D90").Value = 12380
$ws.Range("E90").Value = 1020
$ws.Range("G90").Value = 0
$ws.Range("H90").Value = 320

# Rows 132/133: swap Somalia/Lituania labels similarly
$ws.Range("A132").Value = "Lituania"
$ws.Range("B132").Value = 3397
$ws.Range("C132").Value = 11
$ws.Range("D132").Value = 2094
$ws.Range("E132").Value = 1216
$ws.Range("H132").Value = 87

$ws.Range("A133").Value = "Somalia"
$ws.Range("B133").Value = 3389
$ws.Range("C133").Value = 0
$ws.Range("D133").Value = 2803
$ws.Range("E133").Value = 488
$ws.Range("H133").Value = 98

# Estonia (row 143)
$ws.Range("B143").Value = 2722
$ws.Range("C143").Value = 25
$ws.Range("D143").Value = 2286
$ws.Range("E143").Value = 372
